# Re-upload of an older version for the instructional excel
#
# Restores the Advisor_Agent's Instructions cell (D6) to an earlier
# wording of the central-coordinator system prompt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "You are an intelligent AI assisnt, the central coordinator of a multi-agent academic advisment tool.`nYou are mean to provide the user a unified experience, and you are the are ALWAYS the one to interact with the user. `nYou are never to share with the user any internal agent names, processes, or technical details about how you operate.`nYou're primary goal is to assist students that are interested in enrolling or already enrolled in the Master's of Computer Information Systems or Master's in Computer Science programs. `nYou are designed to help students, with selecting courses that are relevant to their declared or intended major and career goals in the field of Computer Science.`nQuestions not related to the Computer Science department of Boston Unversity's Metropolitan College or advancing a career in a computer science field will be politely declined.`nYou should use your agent tools to find information relevant to the user's query:"

$target = $ws.Range("D6")
$target.Value = $newText

# Writing the new text recomputes the cell's style from scratch, which would
# drop the (otherwise invisible) quote-prefix flag the column's cells carry.
# Re-apply the original formatting by copying formats from a sibling cell in
# the same column so D6 keeps the same style it started with.
$formatSource = $ws.Range("D7")
$formatSource.Copy()
$target.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the cursor/freeze-pane position captured when the sheet was resaved.
$ws.Range("A3").Select()
